$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = '  -3.75%  '
$ws.Range("D2").Value = '60.623.45'

# Row 3
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("D3").Value = '3.348.77'

# Row 4
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.84'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("E6").Value = '  +4.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.92'
$ws.Range("D6").Style = "Normal"

# Row 7
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").Value = '  -1.71%  '
$ws.Range("D8").Value = '3.349.47'

# Row 9
$ws.Range("E9").Value = '  +0.03%  '

# Row 10
$ws.Range("E10").Value = '  +2.08%  '

# Row 11
$ws.Range("E11").Value = '  +0.37%  '

# Row 12
$ws.Range("E12").Value = '  +2.28%  '

# Row 13
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D13").Value = '3.917.63'

# Row 14
$ws.Range("E14").Value = '  +0.67%  '

# Row 15
$ws.Range("E15").Value = '  -0.41%  '

# Row 16
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("D16").Value = '3.346.38'

# Row 17
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.98'
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("D18").Value = '60.741.63'

# Row 19
$ws.Range("E19").Value = '  +5.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.91'
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = '  +2.05%  '

# Row 21
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.26'
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '374.85'
$ws.Range("D22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = '  +0.80%  '

# Row 24
$ws.Range("B24").Value = 'WrappedeETH'
$ws.Range("C24").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("D24").Value = '3.480.42'

# Row 25
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = '  -3.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.13'
$ws.Range("D26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = '  +6.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000115'
$ws.Range("D27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = '  +20.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.68'
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("E29").Value = '  +9.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.67'
$ws.Range("D29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.998'
$ws.Range("D30").Style = "Normal"

# Row 31
$ws.Range("E31").Value = '  +2.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.08'
$ws.Range("D31").Style = "Normal"

# Row 32
$ws.Range("E32").Value = '  +0.00%  '

# Row 33
$ws.Range("E33").Value = '  +0.92%  '

# Row 34
$ws.Range("E34").Value = '  -0.07%  '

# Row 35
$ws.Range("E35").Value = '  -1.87%  '
$ws.Range("D35").Value = '3.377.40'

# Row 36
$ws.Range("E36").Value = '  +1.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.12'
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = '  +3.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.47'
$ws.Range("D37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = '  +3.15%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.95'
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = '  +3.15%  '

# Row 40
$ws.Range("E40").Value = '  -0.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '162.67'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = '  +2.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0780'
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = '  -0.05%  '

# Row 43
$ws.Range("E43").Value = '  +12.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.23'
$ws.Range("D43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = '  +3.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.40'
$ws.Range("D44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = '  -0.72%  '

# Row 46
$ws.Range("E46").Value = '  -3.72%  '

# Row 47
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E47").Value = '  +3.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.45'
$ws.Range("D47").Style = "Normal"

# Row 48
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.59'
$ws.Range("D48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = '  +3.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.93'
$ws.Range("D49").Style = "Normal"

# Row 50
$ws.Range("E50").Value = '  +13.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.98'
$ws.Range("D50").Style = "Normal"

# Row 51
$ws.Range("E51").Value = '  +3.47%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.890'
$ws.Range("D51").Style = "Normal"
